$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: a new price record is inserted at row 5, pushing the
# existing rows 5-25 down to 6-26 (dimension grows from A1:R25 to A1:R26).
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new week's record.
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C5").Value = "Los Lagos"
$ws.Range("D5").Value = 44831
$ws.Range("E5").Value = 10
$ws.Range("F5").Value = 100112035
$ws.Range("G5").Value = "Bruselas (repollito)"
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 90
$ws.Range("K5").Value = 25000
$ws.Range("L5").Value = 25000
$ws.Range("M5").Value = 25000
$ws.Range("N5").Value = "`$/malla 15 kilos"
$ws.Range("O5").Value = "Provincia de Quillota"
$ws.Range("P5").Value = 1667
$ws.Range("Q5").Value = 15
$ws.Range("R5").Value = "Hortaliza"
